$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert E2:E13 formulas into a shared formula group (same formula, re-applied as a block)
$ws.Range("E2:E13").Formula = "=C2/(1.02*0.997)"

# Add a new "Total Sum" row (row 15): label in A15, SUM formula in C15
$ws.Range("A15").Value = "Total Sum"
$ws.Range("C15").Formula = "=SUM(C2:C13)"

# Update the active selection to C16
$ws.Range("C16").Select()
